$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (cohort_year 2020, period_index 5): num_customers 25 -> 26, retention_rate recalculated (26/2654)
$ws.Range("C22").Value = 26
$ws.Range("E22").Value = 0.009796533534287867

# Row 36 (cohort_year 2024, period_index 1): num_customers 113 -> 114, retention_rate recalculated (114/1930)
$ws.Range("C36").Value = 114
$ws.Range("E36").Value = 0.05906735751295337

# Row 37 (cohort_year 2025, period_index 0): num_customers 707 -> 713, cohort_size 707 -> 713 (retention_rate stays 1)
$ws.Range("C37").Value = 713
$ws.Range("D37").Value = 713
